$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 holds a purely numeric-looking id ("113564") that must stay TEXT
# (matches the other CIN/IF cells in the column, e.g. BB125874). Briefly
# force a text number-format so Excel doesn't auto-coerce it to a number,
# then clear the format again so no visible style change is left behind.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "113564"
$ws.Range("C2").ClearFormats()

# --- Row 2 ---
$ws.Range("A2").Value = "001/TEST DR/AV"
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("D2").Value = "lala morale"
$ws.Range("E2").Value = "oui"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = "--"
$ws.Range("H2").Value = 24000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 24000
$ws.Range("M2").Value = 24000

# --- Row 3 ---
$ws.Range("A3").Value = "001/TEST DR/AV"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "BB125874"
$ws.Range("D3").Value = "YASSINE TYEST"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = "--"
$ws.Range("H3").Value = 36000
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = "--"
$ws.Range("K3").Value = 5400
$ws.Range("L3").Value = 36000
$ws.Range("M3").Value = 30600

# --- Row 4 (new) ---
$ws.Range("A4").Value = "002/tEST drrr/AV"
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "BB147852"
$ws.Range("D4").Value = "Mustapha Tahiri"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = "--"
$ws.Range("H4").Value = 20000
$ws.Range("I4").Value = 15
$ws.Range("J4").Value = "--"
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 17000

# --- Row 5 (new) ---
$ws.Range("A5").Value = "002/LF/tEST drrr/AV"
$ws.Range("B5").Value = "Logement de fonction"
$ws.Range("C5").Value = "BB169785"
$ws.Range("D5").Value = "Tawfiq mf"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 3000
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 10
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 2700

# --- Row 6 (new) ---
$ws.Range("A6").Value = "002/LF/tEST drrr/AV"
$ws.Range("B6").Value = "Logement de fonction"
$ws.Range("C6").Value = "BB979797"
$ws.Range("D6").Value = "Test mf"
$ws.Range("E6").Value = "non"
$ws.Range("F6").Value = "mensuelle"
$ws.Range("G6").Value = 7000
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 10
$ws.Range("J6").Value = 700
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 6300
